$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregar nueva persona (fila 5): id=4, alex, 21, 9789485, alex@gmail.com, 01/12/2005
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "alex"
$ws.Cells.Item(5, 3).Value = 21
$ws.Cells.Item(5, 4).Value = 9789485
$ws.Cells.Item(5, 5).Value = "alex@gmail.com"

# La fecha de nacimiento se guarda como texto (no como fecha numerica)
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = "01/12/2005"
$ws.Cells.Item(5, 6).Style = "Normal"

# Mover la seleccion activa a G6
[void]$ws.Range("G6").Select()
